# Update to current US develop branch
# - BBNPPTY sheet: the "new combined cycle gas without CCS" ban (years 2032-2050,
#   columns M:AE) is turned off (set back to 0) for row 4.
# - About sheet: the note that referenced "...banned starting in 2028 and new
#   combined cycle gas without CCS is banned starting in 2032." is simplified so
#   that it just ends in the literal year 2028 (the "2032" sentence is removed).

$wb = $excel.ActiveWorkbook

# --- BBNPPTY sheet: turn off the CCGT-without-CCS ban (M4:AE4 -> 0) ---
$wsData = $wb.Worksheets.Item("BBNPPTY")
$wsData.Range("M4:AE4").Value = 0

# Update the view/selection on the BBNPPTY sheet.
$wsData.Range("L4:AE4").Select() | Out-Null

# --- About sheet: update note text (A13 becomes the literal year 2028) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A13").Value = 2028

# Update the view/selection on the About sheet, and leave it as the active tab
# (it was the active tab before the edit as well).
$wsAbout.Select() | Out-Null
$wsAbout.Range("B18").Select() | Out-Null
